$d = $word.ActiveDocument

# 1. Strip "%" signs from CMYK column values in the Greyscale table
$cmyk = @(
    @('74%, 68%, 67%, 87%', '74, 68, 67, 87'),
    @('72%, 66%, 65%, 74%', '72, 66, 65, 74'),
    @('69%, 62%, 61%, 52%', '69, 62, 61, 52'),
    @('64%, 56%, 56%, 32%', '64, 56, 56, 32'),
    @('58%, 49%, 49%, 16%', '58, 49, 49, 16'),
    @('50%, 41%, 41%, 5%', '50, 41, 41, 5'),
    @('40%, 33%, 33%, 1%', '40, 33, 33, 1'),
    @('17%, 13%, 13%, 0%', '17, 13, 13, 0'),
    @('6%, 4%, 5%, 0%', '6, 4, 5, 0'),
    @('0%, 0%, 0%, 0%', '0, 0, 0, 0')
)
foreach ($pair in $cmyk) {
    $d.Content.Find.Execute($pair[0], $true, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null
}

# 2. greyscaleLevel7 CMYK cell: split into 3 runs "28" / ", " / "22, 22, 0"
$t1 = $d.Tables.Item(1)
$cell97 = $t1.Cell(9, 4)
$pRange97 = $cell97.Range.Paragraphs.Item(1).Range
$frag97 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>28</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>22, 22, 0</w:t></w:r></w:p>'
$pRange97.InsertXML($frag97) | Out-Null

# 3. Merge the 3 "greyscaleLevel"/"1"/"0" runs into one run
$d.Content.Find.Execute("greyscaleLevel10", $true, $false, $false, $false, $false, $true, 1, $false, "greyscaleLevel10", 2) | Out-Null

# 4. Insert "System status colour" heading + new status-colour table
$bookmarkPara = $d.Range($t1.Range.End, $t1.Range.End + 1)
$tableXml = '<w:tbl><w:tblPr><w:tblStyle w:val="LightList-Accent2"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="2394"/><w:gridCol w:w="2394"/><w:gridCol w:w="2394"/><w:gridCol w:w="2394"/></w:tblGrid><w:tr><w:trPr><w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="2394" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Name</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2394" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>Hex</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2394" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>RGB</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2394" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>CMYK</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="2394" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>statusSuccess</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2394" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>#</w:t></w:r><w:r><w:t>58af5a</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2394" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>88, 175, 90</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2394" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>69, 6, 87, 0</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="2394" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>statusWarning</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2394" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>#</w:t></w:r><w:r><w:t>ea8236</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2394" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>234, 130, 54</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2394" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>4, 59, 89, 0</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="2394" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>statusDanger</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2394" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>#</w:t></w:r><w:r><w:t>e8404c</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2394" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>232, 64, 76</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2394" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>3, 90, 69, 0</w:t></w:r></w:p></w:tc></w:tr></w:tbl>'
$fragHeading = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>System status colour</w:t></w:r></w:p>'
$fragTail = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$fullFrag = $fragHeading + $tableXml + $fragTail
$bookmarkPara.InsertXML($fullFrag) | Out-Null

# 5. Move lastRenderedPageBreak from "Gutter size" run to "Layout" run
$rngL = $d.Content
$rngL.Find.Execute("Layout") | Out-Null
$rngL.Expand(4) | Out-Null
$fragLayout = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Layout</w:t></w:r></w:p>'
$rngL.InsertXML($fragLayout) | Out-Null

$rngG = $d.Content
$rngG.Find.Execute("Gutter size") | Out-Null
$rngG.Expand(4) | Out-Null
$fragGutter = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Gutter size</w:t></w:r><w:r><w:t>: 1em</w:t></w:r></w:p>'
$rngG.InsertXML($fragGutter) | Out-Null

Write-Host "Done."
